$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.176.75"
$ws.Range("E2").Value = "  +2.43%  "

$ws.Range("D3").Value = "1.587.40"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("E4").Value = "  +1.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.51%  "

$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("D12").Value = "1.814.64"
$ws.Range("E12").Value = "  +1.44%  "

$ws.Range("D13").Value = "1.582.23"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.529"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").Value = "28.234.69"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("E17").Value = "  +1.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "

$ws.Range("D19").Value = "0.0₃0706"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  -0.62%  "

$ws.Range("D34").Value = "1.398.30"
$ws.Range("E34").Value = "  -4.35%  "

$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("E36").Value = "  -7.34%  "

$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.539"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("E44").Value = "  -1.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.980"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").Value = "1.724.43"
$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("E50").Value = "  +1.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0521"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
